$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for column G (header "K"), rows 2-20, replacing old Strike# values with K values.
$values = @{
    2  = 2
    3  = 5
    4  = 0
    5  = 8
    6  = 6
    7  = 6
    8  = 1
    9  = 4
    10 = 6
    11 = 3
    12 = 5
    13 = 6
    14 = 4
    15 = 4
    16 = 5
    17 = 2
    18 = 3
    19 = 2
    20 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
